$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need an explicit
# Text number format first, otherwise Excel would silently convert
# the string into a float and drop significant trailing zeros
# (e.g. "703.70" -> 703.7, "84.00" -> 84).
$textCells = 'D5', 'D6', 'D9', 'D10', 'D11', 'D13', 'D14', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D40', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '70.905.68'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '3.801.98'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '703.70'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').Value = '170.03'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('D7').Value = '3.800.75'
$ws.Range('E7').Value = '  -1.73%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').Value = '7.59'
$ws.Range('E11').Value = '  +5.37%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  -3.73%  '
$ws.Range('D14').Value = '35.76'
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').Value = '4.442.29'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').Value = '3.832.16'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '70.838.70'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '17.37'
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.114'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '7.10'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = '497.90'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = '10.66'
$ws.Range('E22').Value = '  -4.24%  '
$ws.Range('D23').Value = '0.721'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '84.00'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('E25').Value = '  -5.84%  '
$ws.Range('D26').Value = '3.950.31'
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').Value = '12.02'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('D28').Value = '10.28'
$ws.Range('E28').Value = '  -4.27%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  -6.47%  '
$ws.Range('D31').Value = '3.03'
$ws.Range('E31').Value = '  -4.28%  '
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').Value = '7.32'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').Value = '28.98'
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('E35').Value = '  -4.10%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.766.44'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '9.02'
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').Value = '0.991'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('E39').Value = '  -3.93%  '
$ws.Range('D40').Value = '2.37'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '3.23'
$ws.Range('E44').Value = '  -6.17%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '166.98'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').Value = '0.000313'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').Value = '49.01'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = '416.40'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '8.58'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  -3.57%  '
